$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 28-32 (the five "-inf" rows); rows below shift up to close the gap.
$ws.Range("A28:B32").EntireRow.Delete()
